$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2531.5
$ws.Range("I96").Value = 526
$ws.Range("J96").Value = 3200
$ws.Range("K96").Value = 1578
$ws.Range("L96").Value = 9600
$ws.Range("M96").Value = -205
$ws.Range("N96").Value = -12346

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2647.6667
$ws.Range("I111").Value = 10029
$ws.Range("J111").Value = 1725
$ws.Range("K111").Value = 30087
$ws.Range("L111").Value = 5175
$ws.Range("M111").Value = -27020
$ws.Range("N111").Value = -11309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6424217
$ws.Range("I32").Value = 7474856
$ws.Range("K32").Value = 7474856
$ws.Range("M32").Value = -7474569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 18521566
$ws.Range("I74").Value = 2216.3125
$ws.Range("J74").Value = 45458800
$ws.Range("K74").Value = 2216.3125
$ws.Range("L74").Value = 45458800
$ws.Range("M74").Value = -1342.3125
$ws.Range("N74").Value = -45460548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 18521566
$ws.Range("I77").Value = 2216.3125
$ws.Range("J77").Value = 45458800
$ws.Range("K77").Value = 11081.5625
$ws.Range("L77").Value = 227294000
$ws.Range("M77").Value = -6713.5625
$ws.Range("N77").Value = -227302736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 85157.836
$ws.Range("I122").Value = 144170.58
$ws.Range("J122").Value = 2540
$ws.Range("K122").Value = 432511.74
$ws.Range("L122").Value = 7620
$ws.Range("M122").Value = -430061.74
$ws.Range("N122").Value = -12520

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 78399.766
$ws.Range("I107").Value = 112377.445
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 112377.445
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = -110457.445
$ws.Range("N107").Value = -5790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 68014
$ws.Range("J100").Value = 68014
$ws.Range("L100").Value = 68014
$ws.Range("N100").Value = -70178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25643570
$ws.Range("I132").Value = 33335648
$ws.Range("J132").Value = 15154372
$ws.Range("K132").Value = 100006944
$ws.Range("L132").Value = 45463116
$ws.Range("M132").Value = -100004414
$ws.Range("N132").Value = -45468176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 960.9756
$ws.Range("I5").Value = 748.9048
$ws.Range("J5").Value = 1183.65
$ws.Range("K5").Value = 2246.7144
$ws.Range("L5").Value = 3550.95
$ws.Range("M5").Value = -2134.7144
$ws.Range("N5").Value = -3774.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 22222450
$ws.Range("I50").Value = 199.25
$ws.Range("J50").Value = 111111460
$ws.Range("K50").Value = 597.75
$ws.Range("L50").Value = 333334380
$ws.Range("M50").Value = -116.75
$ws.Range("N50").Value = -333335342

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 22222450
$ws.Range("I53").Value = 199.25
$ws.Range("J53").Value = 111111460
$ws.Range("K53").Value = 597.75
$ws.Range("L53").Value = 333334380
$ws.Range("M53").Value = -116.75
$ws.Range("N53").Value = -333335342

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4500
$ws.Range("I87").Value = 2000
$ws.Range("K87").Value = 6000
$ws.Range("M87").Value = -4752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 4500
$ws.Range("I90").Value = 2000
$ws.Range("K90").Value = 18000
$ws.Range("M90").Value = -11760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2825.804
$ws.Range("I122").Value = 460.94736
$ws.Range("K122").Value = 4148.52624
$ws.Range("M122").Value = -1698.52624

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2605.5454
$ws.Range("J126").Value = 2605.5454
$ws.Range("L126").Value = 7816.6362
$ws.Range("N126").Value = -17696.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3453.6042
$ws.Range("I131").Value = 592.8889
$ws.Range("J131").Value = 4113.769
$ws.Range("K131").Value = 1778.6667
$ws.Range("L131").Value = 12341.307
$ws.Range("M131").Value = 3261.3333
$ws.Range("N131").Value = -22421.307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 960.9756
$ws.Range("I135").Value = 748.9048
$ws.Range("J135").Value = 1183.65
$ws.Range("K135").Value = 6740.1432
$ws.Range("L135").Value = 10652.85
$ws.Range("M135").Value = -4205.1432
$ws.Range("N135").Value = -15722.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1543
$ws.Range("I102").Value = 1647.4445
$ws.Range("J102").Value = 1229.6666
$ws.Range("K102").Value = 1647.4445
$ws.Range("L102").Value = 1229.6666
$ws.Range("M102").Value = -25.44450000000006
$ws.Range("N102").Value = -4473.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 37043844
$ws.Range("I132").Value = 83344280
$ws.Range("J132").Value = 3495.8
$ws.Range("K132").Value = 250032840
$ws.Range("L132").Value = 10487.4
$ws.Range("M132").Value = -250030310
$ws.Range("N132").Value = -15547.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3580.4
$ws.Range("I40").Value = 3187.7144
$ws.Range("J40").Value = 4496.6665
$ws.Range("K40").Value = 3187.7144
$ws.Range("L40").Value = 4496.6665
$ws.Range("M40").Value = -3051.7144
$ws.Range("N40").Value = -4768.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1811.2778
$ws.Range("I68").Value = 1400
$ws.Range("J68").Value = 1862.6875
$ws.Range("K68").Value = 1400
$ws.Range("L68").Value = 1862.6875
$ws.Range("N68").Value = -3360.6875
$ws.Range("M68").Value = -651

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1811.2778
$ws.Range("I71").Value = 1400
$ws.Range("J71").Value = 1862.6875
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 9313.4375
$ws.Range("N71").Value = -16801.4375
$ws.Range("M71").Value = -3256

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 103381.555
$ws.Range("I100").Value = 131672.86
$ws.Range("J100").Value = 4362
$ws.Range("K100").Value = 131672.86
$ws.Range("L100").Value = 4362
$ws.Range("M100").Value = -131131.86
$ws.Range("N100").Value = -5444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2847
$ws.Range("I132").Value = 2610.5557
$ws.Range("J132").Value = 3556.3333
$ws.Range("K132").Value = 7831.6671
$ws.Range("L132").Value = 10668.9999
$ws.Range("M132").Value = -5301.6671
$ws.Range("N132").Value = -15728.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 85000
$ws.Range("J64").Value = 85000
$ws.Range("L64").Value = 85000
$ws.Range("N64").Value = -85496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 85000
$ws.Range("J67").Value = 85000
$ws.Range("L67").Value = 85000
$ws.Range("N67").Value = -86716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 35691
$ws.Range("J92").Value = 35691
$ws.Range("L92").Value = 35691
$ws.Range("N92").Value = -40683

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 44594.5
$ws.Range("J93").Value = 44594.5
$ws.Range("L93").Value = 44594.5
$ws.Range("N93").Value = -49586.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2449.4614
$ws.Range("I122").Value = 2449.4614
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7348.3842
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4898.3842
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5611344.5
$ws.Range("I132").Value = 2470.9429
$ws.Range("J132").Value = 17159024
$ws.Range("K132").Value = 7412.8287
$ws.Range("L132").Value = 51477072
$ws.Range("M132").Value = -4882.8287
$ws.Range("N132").Value = -51482132
